# Atualização de bases das ligas, do dia: 17-02-2024 às 22:47
#
# This refreshes the "Peru Liga 1" odds feed for matches already present
# in the sheet (rows 324-326, 369-372) with newer data pulled from the
# source, and drops three rows (373-375) that no longer exist upstream
# (the whole trailing block shifts up and the tail is removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 324 (match 7302200 -> 7302796 data) ---
$ws.Range("B324").Value = 7302796
$ws.Range("F324").Value = "Sport Huancayo"
$ws.Range("G324").Value = "Sport Boys"
$ws.Range("H324").Value = 1
$ws.Range("I324").Value = 0
$ws.Range("K324").Value = 1.727
$ws.Range("L324").Value = 3.75
$ws.Range("M324").Value = 4.333
$ws.Range("N324").Value = 1.25
$ws.Range("O324").Value = 5.25
$ws.Range("P324").Value = 10
$ws.Range("Q324").Value = -1.75
$ws.Range("R324").Value = 1.925
$ws.Range("S324").Value = 1.875
$ws.Range("T324").Value = 3
$ws.Range("U324").Value = 1.875
$ws.Range("V324").Value = 1.925
$ws.Range("W324").Value = 0.25
$ws.Range("Z324").Value = -1
$ws.Range("AA324").Value = 0.875
$ws.Range("AB324").Value = -1
$ws.Range("AC324").Value = 0.925

# --- Row 325 (match 7302796 -> 7302795 data) ---
$ws.Range("B325").Value = 7302795
$ws.Range("F325").Value = "Unin Comercio"
$ws.Range("G325").Value = "Deportivo Garcilaso"
$ws.Range("I325").Value = 2
$ws.Range("J325").Value = "A"
$ws.Range("K325").Value = 2.25
$ws.Range("L325").Value = 3.3
$ws.Range("M325").Value = 2.7
$ws.Range("N325").Value = 1.75
$ws.Range("O325").Value = 3.6
$ws.Range("P325").Value = 4
$ws.Range("Q325").Value = -0.5
$ws.Range("R325").Value = 1.8
$ws.Range("S325").Value = 2
$ws.Range("T325").Value = 2.75
$ws.Range("U325").Value = 1.825
$ws.Range("V325").Value = 1.975
$ws.Range("W325").Value = -1
$ws.Range("X325").Value = -1
$ws.Range("Y325").Value = 3
$ws.Range("AA325").Value = 1
$ws.Range("AB325").Value = 0.4125
$ws.Range("AC325").Value = -0.5

# --- Row 326 (match 7302795 -> 7302200 data) ---
$ws.Range("B326").Value = 7302200
$ws.Range("F326").Value = "Carlos Manucci"
$ws.Range("G326").Value = "Deportivo Binacional"
$ws.Range("H326").Value = 3
$ws.Range("J326").Value = "H"
$ws.Range("K326").Value = 2
$ws.Range("L326").Value = 3.2
$ws.Range("M326").Value = 3.75
$ws.Range("O326").Value = 3.4
$ws.Range("P326").Value = 4.333
$ws.Range("R326").Value = 1.85
$ws.Range("S326").Value = 1.95
$ws.Range("T326").Value = 2.5
$ws.Range("U326").Value = 1.85
$ws.Range("V326").Value = 1.95
$ws.Range("W326").Value = 0.75
$ws.Range("Y326").Value = -1
$ws.Range("Z326").Value = 0.8500000000000001
$ws.Range("AA326").Value = -1
$ws.Range("AB326").Value = 0.8500000000000001
$ws.Range("AC326").Value = -1

# --- Row 369 ---
$ws.Range("B369").Value = 7768163
$ws.Range("E369").Value = 45340.66666666666
$ws.Range("F369").Value = "Union Comercio"
$ws.Range("G369").Value = "Alianza Lima"
$ws.Range("K369").Value = 3.5
$ws.Range("L369").Value = 3.3
$ws.Range("M369").Value = 2
$ws.Range("N369").Value = 5
$ws.Range("O369").Value = 3.6
$ws.Range("P369").Value = 1.65
$ws.Range("Q369").Value = 0.75
$ws.Range("R369").Value = 2
$ws.Range("S369").Value = 1.85
$ws.Range("T369").Value = 2.25
$ws.Range("U369").Value = 1.85
$ws.Range("V369").Value = 2

# --- Row 370 ---
$ws.Range("B370").Value = 7768164
$ws.Range("E370").Value = 45340.79166666666
$ws.Range("F370").Value = "Cienciano"
$ws.Range("G370").Value = "Sport Boys"
$ws.Range("K370").Value = 1.444
$ws.Range("L370").Value = 4
$ws.Range("M370").Value = 6.5
$ws.Range("N370").Value = 1.4
$ws.Range("O370").Value = 4
$ws.Range("P370").Value = 7.5
$ws.Range("Q370").Value = -1.25
$ws.Range("R370").Value = 2.05
$ws.Range("S370").Value = 1.8
$ws.Range("U370").Value = 2.05
$ws.Range("V370").Value = 1.8

# --- Row 371 ---
$ws.Range("B371").Value = 7768165
$ws.Range("E371").Value = 45340.89583333334
$ws.Range("F371").Value = "Cesar Vallejo"
$ws.Range("G371").Value = "Sport Huancayo"
$ws.Range("K371").Value = 1.833
$ws.Range("L371").Value = 3.6
$ws.Range("M371").Value = 3.75
$ws.Range("N371").Value = 1.4
$ws.Range("O371").Value = 4.75
$ws.Range("P371").Value = 7
$ws.Range("Q371").Value = -1.25
$ws.Range("R371").Value = 1.95
$ws.Range("S371").Value = 1.9
$ws.Range("T371").Value = 2.75
$ws.Range("U371").Value = 1.975
$ws.Range("V371").Value = 1.875

# --- Row 372 ---
$ws.Range("B372").Value = 7768166
$ws.Range("E372").Value = 45341.70833333334
$ws.Range("F372").Value = "Comerciantes Unidos"
$ws.Range("G372").Value = "Atletico Grau"
$ws.Range("K372").Value = 2
$ws.Range("M372").Value = 3.5
$ws.Range("N372").Value = 1.666
$ws.Range("O372").Value = 3.5
$ws.Range("P372").Value = 4.75
$ws.Range("Q372").Value = -0.75
$ws.Range("T372").Value = 2.5
$ws.Range("U372").Value = 2
$ws.Range("V372").Value = 1.85

# --- Drop the trailing rows that no longer exist upstream ---
$ws.Rows("373:375").Delete()
